# support palette format .pal
# Update existing "prompt to save" task (row 14 on the Active sheet) to note
# that the tracking needs to be added to the MasterImage class, and add four
# new documentation-related todo items. Also bump the tracked "Max Id" value
# on the Config sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Active")

# 1. Amend the existing task text in B14.
$ws.Range("B14").Value = "prompt to save if image has changed since last save`n- on closing program`n- on opening new image`nAdd this tracking to MasterImage class"

# 2. Insert 4 new rows right after row 15 (before the old row 16) to hold the
#    new documentation tasks, pushing everything below down by 4 rows.
$ws.Rows("16:19").Insert()

# 3. Fill in the new rows with the new todo items.
$ws.Range("A16").Value = 76
$ws.Range("B16").Value = "full documentation of Perpetual Paint"
$ws.Range("C16").Value = "Todo"
$ws.Range("D16").Value = "Task"
$ws.Range("E16").Value = "'8/24/2018"

$ws.Range("A17").Value = 79
$ws.Range("B17").Value = "documentation: include request for sample palette files for the specific formats/color spaces I can't verify because I don't have a test file to load`n- maybe in the actual error messages from the library, too`n- like, send me the file you are trying to read so I can add support for it"
$ws.Range("C17").Value = "Todo"
$ws.Range("D17").Value = "Task"
$ws.Range("E17").Value = "'8/24/2018"

$ws.Range("A18").Value = 77
$ws.Range("B18").Value = "full documentation of support library Colors"
$ws.Range("C18").Value = "Todo"
$ws.Range("D18").Value = "Task"
$ws.Range("E18").Value = "'8/24/2018"

$ws.Range("A19").Value = 78
$ws.Range("B19").Value = "full documentation of support library GUI"
$ws.Range("C19").Value = "Todo"
$ws.Range("D19").Value = "Task"
$ws.Range("E19").Value = "'8/24/2018"

# Reset the number formatting so these new cells look like the rest of the
# table (plain text, default style) rather than picking up a "text" format
# from the apostrophe-prefixed date entry above.
$ws.Range("A16:E19").Style = "Normal"

# 4. Update the tracked Max Id on the Config sheet to reflect the new highest id (79).
$cfg = $wb.Worksheets.Item("Config")
$cfg.Range("F2").Value = 79
